# faculty_info.xlsx: normalize "研究方向" (research-area) delimiters from
# semicolons to commas, fix a missing space after a comma, and update the
# last-saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Jianye Hao: "; " -> ", "
$ws.Range("D2").Value = "Artificial intelligence, Reinforcement Learning, Embodied AI, LLM Agent"

# Row 3 - Yan Zheng: "; " -> ", " (keep the preserved trailing space)
$ws.Range("D3").Value = "Multiagent Systems, Deep Reinforcement Learning, Evolutionary Algorithm "

# Row 4 - Hongyao Tang: "; " -> ", "
$ws.Range("D4").Value = "Deep Reinforcement Learning, Representation Learning"

# Row 6 - Tianpei Yang: "; " -> ", "
$ws.Range("D6").Value = "Reinforcement Learning, Transfer Learning, Multiagent Learning"

# Row 7 - Spike/Fei Ni: add missing space after comma in the position cell ...
$ws.Range("C7").Value = "Postdoctoral fellow, Imperial College London"
# ... and "; " -> ", " in the research-area cell
$ws.Range("D7").Value = "Model based RL, Diffusion for RL, LLM"

# Re-apply the existing wrap formatting to C7 so it resolves to the same
# (deduplicated) cell style used elsewhere, matching the saved workbook.
$ws.Range("C7").WrapText = $true

# Last-saved cursor position moved from C7 to E8.
$ws.Range("E8").Select()
